$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K1").EntireColumn.Insert()
